$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033083325388348
$ws.Range("D2").Value = 1.041284524813878
$ws.Range("E2").Value = 1.041860442800671
$ws.Range("F2").Value = 1.051427521299317
$ws.Range("I2").Value = 1.02733236633473
$ws.Range("J2").Value = 1.038209999563212
$ws.Range("K2").Value = 1.044064474674538
$ws.Range("L2").Value = 1.044638762988611
$ws.Range("M2").Value = 1.0541790527236
$ws.Range("N2").Value = 1.016646001481358

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034122276959878
$ws.Range("D3").Value = 1.042205332277166
$ws.Range("E3").Value = 1.042773807360731
$ws.Range("F3").Value = 1.052402219717811
$ws.Range("I3").Value = 1.027293588095579
$ws.Range("J3").Value = 1.038891062454806
$ws.Range("K3").Value = 1.044795589009561
$ws.Range("L3").Value = 1.045362572105777
$ws.Range("M3").Value = 1.054965978487567
$ws.Range("N3").Value = 1.016877661532564

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034795240108966
$ws.Range("D4").Value = 1.04280198210633
$ws.Range("E4").Value = 1.04336564257405
$ws.Range("F4").Value = 1.053033470350272
$ws.Range("I4").Value = 1.027266154843311
$ws.Range("J4").Value = 1.039331876145136
$ws.Range("K4").Value = 1.045268866121904
$ws.Range("L4").Value = 1.045831121922233
$ws.Range("M4").Value = 1.055475108838508
$ws.Range("N4").Value = 1.017027454892993

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035078318932007
$ws.Range("D5").Value = 1.043053010270241
$ws.Range("E5").Value = 1.043614646992806
$ws.Range("F5").Value = 1.053298979754213
$ws.Range("I5").Value = 1.027254060371935
$ws.Range("J5").Value = 1.039517222266748
$ws.Range("K5").Value = 1.045467878103157
$ws.Range("L5").Value = 1.046028146452947
$ws.Range("M5").Value = 1.055689130845623
$ws.Range("N5").Value = 1.017090402231243

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035125858778473
$ws.Range("D6").Value = 1.043095170498022
$ws.Range("E6").Value = 1.043656467442851
$ws.Range("F6").Value = 1.053343567630933
$ws.Range("I6").Value = 1.027251996697445
$ws.Range("J6").Value = 1.039548344317778
$ws.Range("K6").Value = 1.045501295777294
$ws.Range("L6").Value = 1.046061230419332
$ws.Range("M6").Value = 1.055725065084636
$ws.Range("N6").Value = 1.017100969843193

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034799021973953
$ws.Range("D7").Value = 1.04280533558596
$ws.Range("E7").Value = 1.043368969011076
$ws.Range("F7").Value = 1.053037017584413
$ws.Range("I7").Value = 1.027265995444184
$ws.Range("J7").Value = 1.039334352639881
$ws.Range("K7").Value = 1.045271525149064
$ws.Range("L7").Value = 1.045833754393377
$ws.Range("M7").Value = 1.055477968677048
$ws.Range("N7").Value = 1.017028296100212

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033434300317005
$ws.Range("D8").Value = 1.041595544335294
$ws.Range("E8").Value = 1.042168946644643
$ws.Range("F8").Value = 1.051756809334419
$ws.Range("I8").Value = 1.02731974516229
$ws.Range("J8").Value = 1.038440142538302
$ws.Range("K8").Value = 1.044311516780978
$ws.Range("L8").Value = 1.04488333632153
$ws.Range("M8").Value = 1.054445010311114
$ws.Range("N8").Value = 1.016724313843911

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031034810697724
$ws.Range("D9").Value = 1.039470112173963
$ws.Range("E9").Value = 1.040060740659577
$ws.Range("F9").Value = 1.049505234159502
$ws.Range("I9").Value = 1.02739657497971
$ws.Range("J9").Value = 1.036865387506529
$ws.Range("K9").Value = 1.042621411398038
$ws.Range("L9").Value = 1.04321012843336
$ws.Range("M9").Value = 1.052624376494518
$ws.Range("N9").Value = 1.016187861319516

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029438761604544
$ws.Range("D10").Value = 1.03805750917663
$ws.Range("E10").Value = 1.038659631503711
$ws.Range("F10").Value = 1.04800716534549
$ws.Range("I10").Value = 1.027435832536215
$ws.Range("J10").Value = 1.035816240430935
$ws.Range("K10").Value = 1.041495772584287
$ws.Range("L10").Value = 1.042095754364369
$ws.Range("M10").Value = 1.051410405220917
$ws.Range("N10").Value = 1.015829712085464

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02874851700642
$ws.Range("D11").Value = 1.03744688217315
$ws.Range("E11").Value = 1.03805398410695
$ws.Range("F11").Value = 1.047359208568668
$ws.Range("I11").Value = 1.027450007092135
$ws.Range("J11").Value = 1.03536212105457
$ws.Range("K11").Value = 1.0410086301701
$ws.Range("L11").Value = 1.041613489865786
$ws.Range("M11").Value = 1.050884706224684
$ws.Range("N11").Value = 1.015674511988107

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028492258309641
$ws.Range("D12").Value = 1.037220225324184
$ws.Range("E12").Value = 1.03782917734006
$ws.Range("F12").Value = 1.047118637555887
$ws.Range("I12").Value = 1.02745484885297
$ws.Range("J12").Value = 1.035193466879897
$ws.Range("K12").Value = 1.04082772470032
$ws.Range("L12").Value = 1.041434396320897
$ws.Range("M12").Value = 1.050689433153749
$ws.Range("N12").Value = 1.015616846207986

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028547220831454
$ws.Range("D13").Value = 1.037268836856856
$ws.Range("E13").Value = 1.037877392003723
$ws.Range("F13").Value = 1.047170235902489
$ws.Range("I13").Value = 1.02745382941801
$ws.Range("J13").Value = 1.035229642585803
$ws.Range("K13").Value = 1.040866527672539
$ws.Range("L13").Value = 1.041472810626521
$ws.Range("M13").Value = 1.050731320100458
$ws.Range("N13").Value = 1.015629216502423

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028727331951169
$ws.Range("D14").Value = 1.037428143434245
$ws.Range("E14").Value = 1.038035398287378
$ws.Range("F14").Value = 1.047339320663072
$ws.Range("I14").Value = 1.02745041594301
$ws.Range("J14").Value = 1.035348179513776
$ws.Range("K14").Value = 1.04099367562797
$ws.Range("L14").Value = 1.041598685096809
$ws.Range("M14").Value = 1.050868564988146
$ws.Range("N14").Value = 1.015669745672845

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028838321446296
$ws.Range("D15").Value = 1.037526318337329
$ws.Range("E15").Value = 1.038132772095009
$ws.Range("F15").Value = 1.047443513844251
$ws.Range("I15").Value = 1.027448256726937
$ws.Range("J15").Value = 1.035421217491254
$ws.Range("K15").Value = 1.041072021126693
$ws.Range("L15").Value = 1.041676245966469
$ws.Range("M15").Value = 1.050953125457865
$ws.Range("N15").Value = 1.015694714716571

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029484588916931
$ws.Range("D16").Value = 1.03809805646853
$ws.Range("E16").Value = 1.03869984837484
$ws.Range("F16").Value = 1.048050183316216
$ws.Range("I16").Value = 1.027434832400708
$ws.Range("J16").Value = 1.035846382443053
$ws.Range("K16").Value = 1.041528108312407
$ws.Range("L16").Value = 1.04212776636726
$ws.Range("M16").Value = 1.051445293358803
$ws.Range("N16").Value = 1.015840009724362

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029890204664883
$ws.Range("D17").Value = 1.038456971713958
$ws.Range("E17").Value = 1.039055840106923
$ws.Range("F17").Value = 1.048430923886227
$ws.Range("I17").Value = 1.027425656265166
$ws.Range("J17").Value = 1.036113122723286
$ws.Range("K17").Value = 1.041814271740535
$ws.Range("L17").Value = 1.042411065240536
$ws.Range("M17").Value = 1.051754007013925
$ws.Range("N17").Value = 1.015931117795056

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030126875954799
$ws.Range("D18").Value = 1.038666421122527
$ws.Range("E18").Value = 1.039263584537949
$ws.Range("F18").Value = 1.048653072504421
$ws.Range("I18").Value = 1.027420031416848
$ws.Range("J18").Value = 1.03626872400021
$ws.Range("K18").Value = 1.041981211694184
$ws.Range("L18").Value = 1.042576334366507
$ws.Range("M18").Value = 1.051934070390959
$ws.Range("N18").Value = 1.015984248075403

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030207588727693
$ws.Range("D19").Value = 1.038737854925088
$ws.Range("E19").Value = 1.039334437030154
$ws.Range("F19").Value = 1.048728831121593
$ws.Range("I19").Value = 1.027418067226576
$ws.Range("J19").Value = 1.036321782738218
$ws.Range("K19").Value = 1.042038138231888
$ws.Range("L19").Value = 1.042632691197749
$ws.Range("M19").Value = 1.051995466622592
$ws.Range("N19").Value = 1.016002362164293

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029846677398914
$ws.Range("D20").Value = 1.038418453131197
$ws.Range("E20").Value = 1.039017635162547
$ws.Range("F20").Value = 1.048390066871247
$ws.Range("I20").Value = 1.027426668962902
$ws.Range("J20").Value = 1.036084502336447
$ws.Range("K20").Value = 1.041783566451774
$ws.Range("L20").Value = 1.042380667289332
$ws.Range("M20").Value = 1.051720885373376
$ws.Range("N20").Value = 1.015921343949963

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028674290154623
$ws.Range("D21").Value = 1.037381227269681
$ws.Range("E21").Value = 1.037988865014684
$ws.Range("F21").Value = 1.047289526396651
$ws.Range("I21").Value = 1.027451432803445
$ws.Range("J21").Value = 1.035313272620049
$ws.Range("K21").Value = 1.040956232587902
$ws.Range("L21").Value = 1.041561617064849
$ws.Range("M21").Value = 1.050828149919698
$ws.Range("N21").Value = 1.015657811324139

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027937909623808
$ws.Range("D22").Value = 1.036729992280887
$ws.Range("E22").Value = 1.037342948996361
$ws.Range("F22").Value = 1.046598204101677
$ws.Range("I22").Value = 1.027464554218061
$ws.Range("J22").Value = 1.034828520558251
$ws.Range("K22").Value = 1.040436292084302
$ws.Range("L22").Value = 1.041046885111655
$ws.Range("M22").Value = 1.050266823025508
$ws.Range("N22").Value = 1.015492016511784

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028328208042741
$ws.Range("D23").Value = 1.037075137655422
$ws.Range("E23").Value = 1.037685274434744
$ws.Range("F23").Value = 1.046964626875541
$ws.Range("I23").Value = 1.027457830072975
$ws.Range("J23").Value = 1.035085482279923
$ws.Range("K23").Value = 1.040711899619378
$ws.Range("L23").Value = 1.041319731461187
$ws.Range("M23").Value = 1.050564395362462
$ws.Range("N23").Value = 1.015579916975261

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029866345241864
$ws.Range("D24").Value = 1.038435857711133
$ws.Range("E24").Value = 1.039034898022004
$ws.Range("F24").Value = 1.048408528184251
$ws.Range("I24").Value = 1.027426212210668
$ws.Range("J24").Value = 1.036097434607322
$ws.Range("K24").Value = 1.041797440770798
$ws.Range("L24").Value = 1.042394402735689
$ws.Range("M24").Value = 1.051735851629684
$ws.Range("N24").Value = 1.015925760365454

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031654502158514
$ws.Range("D25").Value = 1.040018825126043
$ws.Range("E25").Value = 1.040604998631088
$ws.Range("F25").Value = 1.050086800683031
$ws.Range("I25").Value = 1.027378825560503
$ws.Range("J25").Value = 1.037272381243577
$ws.Range("K25").Value = 1.043058154601287
$ws.Range("L25").Value = 1.043642503344677
$ws.Range("M25").Value = 1.053095098291619
$ws.Range("N25").Value = 1.016326639371141
